# #5: insurance, claim, debt, investment done
# Rebuild the "保險" (insurance) sheet (sheet6) with the normalized
# property-report column layout: company, name, owner, property_category,
# category, date, legislator_name, legislator_id, source_file, index.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(6)

$newCols = @("F","G","H","I","J","K")

# ---- Header row (row 1) ------------------------------------------------
$ws.Range("B1").Value = "company"
$ws.Range("C1").Value = "name"
$ws.Range("D1").Value = "owner"
$ws.Range("E1").Value = "property_category"
$ws.Range("F1").Value = "category"
$ws.Range("G1").Value = "date"
$ws.Range("H1").Value = "legislator_name"
$ws.Range("I1").Value = "legislator_id"
$ws.Range("J1").Value = "source_file"
$ws.Range("K1").Value = "index"

# Give the newly added header cells (F1:K1) the same bold / centered /
# bordered look already used by B1:E1.
foreach ($c in $newCols) {
    $cell = $ws.Range("$c`1")
    $cell.Font.Bold = $true
    $cell.HorizontalAlignment = -4108
    $cell.VerticalAlignment = -4160
    $cell.Borders.LineStyle = 1
}

# ---- Data rows -----------------------------------------------------------
# Column A already holds the running index (84..93) from the source data.
$data = @(
    @{ Row=2;  Idx=84; Company="國華人壽";       Name="添福增額終身壽險十五年期";         Owner="楊曜" },
    @{ Row=3;  Idx=85; Company="三商美邦人壽";   Name="世紀理財變額萬能終身壽險二十年期"; Owner="楊曜" },
    @{ Row=4;  Idx=86; Company="三商美邦人壽";   Name="世紀理財變額萬能終身壽險二十年期"; Owner="林君倩" },
    @{ Row=5;  Idx=87; Company="遠雄人壽";       Name="富貴一生终身還本保險十年期";       Owner="林君倩" },
    @{ Row=6;  Idx=88; Company="富邦人壽";       Name="新喜樂養老保險二十年期";           Owner="林君倩" },
    @{ Row=7;  Idx=89; Company="富邦人壽";       Name="如意增額養老保險十五年期";         Owner="林君倩" },
    @{ Row=8;  Idx=90; Company="富邦人壽";       Name="新吉富養老保險二十年期";           Owner="林君倩" },
    @{ Row=9;  Idx=91; Company="蘇黎世國際人壽"; Name="傳愛增額終身壽險十五年期";         Owner="林君倩" },
    @{ Row=10; Idx=92; Company="蘇黎世國際人壽"; Name="傳愛增額終身壽險十五年期";         Owner="林君倩" },
    @{ Row=11; Idx=93; Company="安聯人壽";       Name="吉利長红變額萬能壽險二十年期";     Owner="林君倩" }
)

foreach ($r in $data) {
    $row = $r.Row

    $ws.Range("A$row").Value = $r.Idx
    $ws.Range("B$row").Value = $r.Company
    $ws.Range("C$row").Value = $r.Name
    $ws.Range("D$row").Value = $r.Owner
    $ws.Range("E$row").Value = "insurance"
    $ws.Range("F$row").Value = "normal"

    # Keep the date a plain text value ("2012-04-20"), matching every other
    # sheet in this workbook, instead of letting Excel auto-convert it to a
    # date serial number.
    $dcell = $ws.Range("G$row")
    $dcell.NumberFormat = "@"
    $dcell.Value = "2012-04-20"

    $ws.Range("H$row").Value = "楊曜"
    $ws.Range("I$row").Value = 1759
    $ws.Range("J$row").Value = "tmpcdc61"
    $ws.Range("K$row").Value = $r.Idx
}
